$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The "News API Email" application occupies row 9 (columns A-D: Title,
# Description, URL, Image). Select the whole row, the way a user would
# right before deleting it, then remove it entirely -- this shifts every
# row below it up by one and drops the row's three now-unused shared
# strings automatically.
$ws.Range("A9:XFD9").Select()
$ws.Rows.Item(9).Delete()

# Rebuild the hyperlink list so rId8/rId9/rId10 point at the surviving
# rows' URLs (PDF-invoices, happiness-web-app, weather-forecast-web-app)
# instead of the stale news-api-email target left behind by the row
# shift.
$urls = @(
  "https://github.com/valenpendragon/my-web-todo-app",
  "https://github.com/valenpendragon/todo-list",
  "https://github.com/valenpendragon/blackjack-py27",
  "https://github.com/valenpendragon/convert-table-to-cl-md",
  "https://github.com/valenpendragon/history-skeleton-generator",
  "https://github.com/valenpendragon/web-weather-api",
  "https://github.com/valenpendragon/apod-web-page",
  "https://github.com/valenpendragon/PDF-invoices",
  "https://github.com/valenpendragon/happiness-web-app",
  "https://github.com/valenpendragon/weather-forecast-web-app"
)
$ws.Range("C2:C11").Hyperlinks.Delete()
for ($i = 0; $i -lt $urls.Length; $i++) {
    $row = $i + 2
    $ws.Hyperlinks.Add($ws.Range("C$row"), $urls[$i]) | Out-Null
}
